# Auto commit at 2025-12-03  7:39:20.31
# Append two new daily data rows (row 4 and row 5) to Sheet1, mirroring the
# layout of the existing rows (A=date, B=site name, C=charge kwh,
# D=total charge revenue, E=service-fee revenue, F=order count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - 四方坪站
$ws.Range("A4").Value = 45993
$ws.Range("B4").Value = "四方坪站"
$ws.Range("C4").Value = 7975.17
$ws.Range("D4").Value = 6917.88
$ws.Range("E4").Value = 2647.03
$ws.Range("F4").Value = 352

# Row 5 - 高岭站
$ws.Range("A5").Value = 45993
$ws.Range("B5").Value = "高岭站"
$ws.Range("C5").Value = 5614.48
$ws.Range("D5").Value = 4868.16
$ws.Range("E5").Value = 1442.2
$ws.Range("F5").Value = 211

# Match the new active-cell selection recorded in the workbook view.
[void]$ws.Range("F10").Select()
